# week 50 sheet: log a new activity (row 10), extend totals, and move selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 50")

# New log entry in row 10: 10:45 - 11:30, activity description in column F.
$ws.Range("C10").Value = 0.44791666666666669
$ws.Range("D10").Value = 0.47916666666666669
$ws.Range("F10").Value = "Menu class aangemaakt en eerst button aangemaakt."

# Row 10 becomes a wrapped, taller row like row 8.
$ws.Rows.Item(10).RowHeight = 28.5

# Vertically center the time/duration columns that already use the
# horizontal-center + time-format style, matching rows 8 onward.
$ws.Range("G7:G9").VerticalAlignment = -4108
$ws.Range("C10:D18").VerticalAlignment = -4108
$ws.Range("G11:G17").VerticalAlignment = -4108

# New trailing blank row with the vertical-center-only style.
$ws.Range("C19:D19").VerticalAlignment = -4108

# Selection moves to F13 (where the next entry will go).
$ws.Activate() | Out-Null
$ws.Range("F13").Select() | Out-Null

Write-Output "week 50 sheet updated"
